$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (previously derived from Strike# data) for rows 2-68, column G.
$newK = @{
    2 = 8
    3 = 3
    4 = 5
    5 = 0
    6 = 0
    7 = 4
    8 = 1
    9 = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 3
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 3
    19 = 1
    20 = 2
    21 = 2
    22 = 2
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 0
    38 = 1
    39 = 1
    40 = 2
    41 = 0
    42 = 0
    43 = 2
    44 = 1
    45 = 2
    46 = 0
    47 = 0
    48 = 1
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 1
    57 = 1
    58 = 2
    59 = 5
    60 = 2
    61 = 6
    62 = 2
    63 = 4
    64 = 5
    65 = 8
    66 = 3
    67 = 4
    68 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
